# Updates the cryptos list with latest price/volume data
# (as scraped by the GitHub Actions job on Sun May 28 23:00:08 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Price-column (D) updates whose new values are plain numeric strings --
# Force these cells to Text format first so Excel keeps the exact display
# string (e.g. "156.54") instead of silently coercing it to a Double and
# rounding/reformatting it. The format is reset back to Normal afterwards
# so the cells end up with no explicit style, matching their original state.
$textForceCells = @("D4", "D5", "D8", "D9", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$dPriceUpdates = @{
    "D4" = "1.006"
    "D5" = "316.15"
    "D8" = "0.3845"
    "D9" = "0.07392"
    "D11" = "21.02"
    "D12" = "0.07812"
    "D14" = "5.564"
    "D15" = "6.669"
    "D16" = "92.18"
    "D17" = "1.007"
    "D18" = "0.000008914"
    "D21" = "15.09"
    "D22" = "5.177"
    "D24" = "10.93"
    "D25" = "156.54"
    "D26" = "1.934"
    "D27" = "18.62"
    "D28" = "2.130"
    "D29" = "116.96"
    "D30" = "5.054"
    "D31" = "0.08921"
    "D32" = "3.360"
    "D34" = "0.7855"
    "D35" = "4.687"
    "D36" = "2.785"
    "D37" = "1.135"
    "D38" = "0.02061"
    "D39" = "0.05408"
    "D40" = "0.5611"
    "D42" = "7.153"
    "D43" = "8.642"
    "D44" = "0.1542"
    "D45" = "0.4971"
    "D46" = "10.88"
    "D47" = "107.82"
    "D48" = "1.685"
    "D50" = "69.59"
    "D51" = "0.06142"
}
foreach ($ref in $dPriceUpdates.Keys) {
    $ws.Range($ref).Value = $dPriceUpdates[$ref]
}

foreach ($ref in $textForceCells) {
    $ws.Range($ref).Style = "Normal"
}

# -- All other updates (coin name/link swaps, non-numeric-looking prices,
#    and volume/percentage-change text) can be set directly. --
$otherUpdates = @{
    "D2" = "28.228.15"
    "E2" = "  +3.48%  "
    "D3" = "1.915.38"
    "E3" = "  +2.78%  "
    "E4" = "  -1.53%  "
    "E5" = "  +0.99%  "
    "E6" = "  -1.31%  "
    "E7" = "  +0.91%  "
    "E8" = "  +2.91%  "
    "E9" = "  -0.73%  "
    "E10" = "  +1.75%  "
    "E11" = "  +1.68%  "
    "E12" = "  -0.90%  "
    "D13" = "1.920.97"
    "E13" = "  +3.10%  "
    "E14" = "  +2.38%  "
    "E15" = "  +1.93%  "
    "E16" = "  +2.03%  "
    "E17" = "  -1.48%  "
    "E18" = "  +1.37%  "
    "E19" = "  -1.25%  "
    "D20" = "28.219.86"
    "E20" = "  +3.31%  "
    "E21" = "  +1.85%  "
    "E22" = "  +1.10%  "
    "D23" = "2.151.65"
    "E23" = "  +2.59%  "
    "E24" = "  +2.26%  "
    "B25" = "Monero"
    "C25" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "E25" = "  +1.51%  "
    "B26" = "Toncoin"
    "C26" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "E26" = "  -1.11%  "
    "E27" = "  +0.33%  "
    "E28" = "  +6.13%  "
    "E29" = "  +0.78%  "
    "E30" = "  +1.13%  "
    "E31" = "  +0.02%  "
    "E32" = "  +0.29%  "
    "E33" = "  +4.84%  "
    "E34" = "  +5.35%  "
    "E35" = "  +2.57%  "
    "E36" = "  +4.04%  "
    "E37" = "  +0.78%  "
    "E38" = "  +0.30%  "
    "E40" = "  +4.63%  "
    "E41" = "  +0.74%  "
    "E42" = "  +0.28%  "
    "E43" = "  +2.98%  "
    "E44" = "  +0.37%  "
    "E45" = "  +2.86%  "
    "E46" = "  +2.35%  "
    "E48" = "  +0.94%  "
    "E49" = "  -1.36%  "
    "E50" = "  +4.30%  "
    "E51" = "  +0.98%  "
}
foreach ($ref in $otherUpdates.Keys) {
    $ws.Range($ref).Value = $otherUpdates[$ref]
}

Write-Host "Applied $($dPriceUpdates.Count + $otherUpdates.Count) cell updates"
